$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 113.42857
$ws.Range("I5").Value = 113.42857
$ws.Range("K5").Value = 113.42857
$ws.Range("M5").Value = 1.571430000000007
$ws.Range("H18").Value = 638.6
$ws.Range("I18").Value = 298.25
$ws.Range("K18").Value = 298.25
$ws.Range("M18").Value = -14.25
$ws.Range("H62").Value = 7745.8423
$ws.Range("I62").Value = 7204.8237
$ws.Range("K62").Value = 7204.8237
$ws.Range("M62").Value = -6580.8237
$ws.Range("H65").Value = 7745.8423
$ws.Range("I65").Value = 7204.8237
$ws.Range("K65").Value = 36024.1185
$ws.Range("M65").Value = -32904.1185
$ws.Range("H100").Value = 2410.4092
$ws.Range("J100").Value = 4916.6665
$ws.Range("L100").Value = 4916.6665
$ws.Range("N100").Value = -5998.6665
$ws.Range("H112").Value = 45803.74
$ws.Range("I112").Value = 2699.6
$ws.Range("J112").Value = 57777.11
$ws.Range("K112").Value = 8098.799999999999
$ws.Range("L112").Value = 173331.33
$ws.Range("M112").Value = -6990.799999999999
$ws.Range("N112").Value = -175547.33
$ws.Range("H132").Value = 2192.1765
$ws.Range("I132").Value = 2079.1875
$ws.Range("K132").Value = 6237.5625
$ws.Range("M132").Value = -3707.5625
$ws.Range("H137").Value = 2960.4
$ws.Range("I137").Value = 1929.9
$ws.Range("K137").Value = 5789.700000000001
$ws.Range("M137").Value = -3239.700000000001
$ws.Range("H138").Value = 2880.29
$ws.Range("I138").Value = 2403.4
$ws.Range("J138").Value = 2933.2778
$ws.Range("K138").Value = 7210.200000000001
$ws.Range("L138").Value = 8799.8334
$ws.Range("M138").Value = -2070.200000000001
$ws.Range("N138").Value = -19079.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1401.037
$ws.Range("I2").Value = 840.05
$ws.Range("J2").Value = 3003.8572
$ws.Range("K2").Value = 840.05
$ws.Range("L2").Value = 3003.8572
$ws.Range("M2").Value = -727.05
$ws.Range("N2").Value = -3229.8572
$ws.Range("H5").Value = 219.5
$ws.Range("I5").Value = 232.77777
$ws.Range("K5").Value = 232.77777
$ws.Range("M5").Value = -120.77777
$ws.Range("H45").Value = 4057.9119
$ws.Range("I45").Value = 3682.4375
$ws.Range("J45").Value = 4391.6665
$ws.Range("K45").Value = 3682.4375
$ws.Range("L45").Value = 4391.6665
$ws.Range("M45").Value = -3305.4375
$ws.Range("N45").Value = -5145.6665
$ws.Range("H61").Value = 7644.6895
$ws.Range("I61").Value = 9680.388999999999
$ws.Range("K61").Value = 9680.388999999999
$ws.Range("M61").Value = -9468.388999999999
$ws.Range("H97").Value = 1496.5
$ws.Range("I97").Value = 1187.3636
$ws.Range("K97").Value = 1187.3636
$ws.Range("M97").Value = -691.3635999999999
$ws.Range("H116").Value = 1401.037
$ws.Range("I116").Value = 840.05
$ws.Range("J116").Value = 3003.8572
$ws.Range("K116").Value = 840.05
$ws.Range("L116").Value = 3003.8572
$ws.Range("M116").Value = 1453.95
$ws.Range("N116").Value = -7591.8572
$ws.Range("H136").Value = 7644.6895
$ws.Range("I136").Value = 9680.388999999999
$ws.Range("K136").Value = 29041.167
$ws.Range("M136").Value = -26491.167

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1401.037
$ws.Range("I3").Value = 840.05
$ws.Range("J3").Value = 3003.8572
$ws.Range("K3").Value = 840.05
$ws.Range("L3").Value = 3003.8572
$ws.Range("M3").Value = -726.05
$ws.Range("N3").Value = -3231.8572
$ws.Range("H4").Value = 219.5
$ws.Range("I4").Value = 232.77777
$ws.Range("K4").Value = 232.77777
$ws.Range("M4").Value = -117.77777
$ws.Range("H105").Value = 1869.0526
$ws.Range("I105").Value = 1647.8823
$ws.Range("J105").Value = 3749
$ws.Range("K105").Value = 1647.8823
$ws.Range("L105").Value = 3749
$ws.Range("M105").Value = 99.11770000000001
$ws.Range("N105").Value = -7243
$ws.Range("H107").Value = 1715.05
$ws.Range("I107").Value = 1721.0294
$ws.Range("K107").Value = 1721.0294
$ws.Range("M107").Value = 198.9706000000001
$ws.Range("H132").Value = 92500
$ws.Range("J132").Value = 92500
$ws.Range("L132").Value = 92500
$ws.Range("N132").Value = -102620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 28.23077
$ws.Range("I7").Value = 29.818182
$ws.Range("J7").Value = 19.5
$ws.Range("K7").Value = 29.818182
$ws.Range("L7").Value = 19.5
$ws.Range("M7").Value = 83.18181799999999
$ws.Range("N7").Value = -245.5
$ws.Range("H19").Value = 1088
$ws.Range("I19").Value = 643.4
$ws.Range("K19").Value = 643.4
$ws.Range("M19").Value = -473.4
$ws.Range("H22").Value = 2774.0833
$ws.Range("I22").Value = 2889.9092
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 2889.9092
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -2539.9092
$ws.Range("N22").Value = -2200
$ws.Range("H24").Value = 1088
$ws.Range("I24").Value = 643.4
$ws.Range("K24").Value = 643.4
$ws.Range("M24").Value = -473.4
$ws.Range("H31").Value = 4651.795
$ws.Range("I31").Value = 2207.2307
$ws.Range("J31").Value = 5874.077
$ws.Range("K31").Value = 2207.2307
$ws.Range("L31").Value = 5874.077
$ws.Range("M31").Value = -1912.2307
$ws.Range("N31").Value = -6464.077
$ws.Range("H34").Value = 4651.795
$ws.Range("I34").Value = 2207.2307
$ws.Range("J34").Value = 5874.077
$ws.Range("K34").Value = 2207.2307
$ws.Range("L34").Value = 5874.077
$ws.Range("M34").Value = -2005.2307
$ws.Range("N34").Value = -6278.077
$ws.Range("H107").Value = 1878.8948
$ws.Range("I107").Value = 900.3333
$ws.Range("J107").Value = 2330.5386
$ws.Range("K107").Value = 900.3333
$ws.Range("L107").Value = 2330.5386
$ws.Range("M107").Value = 1019.6667
$ws.Range("N107").Value = -6170.5386
$ws.Range("H132").Value = 3334.625
$ws.Range("I132").Value = 3334.625
$ws.Range("K132").Value = 10003.875
$ws.Range("M132").Value = -7473.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2288.5715
$ws.Range("J21").Value = 2473.3333
$ws.Range("L21").Value = 7419.999899999999
$ws.Range("N21").Value = -7765.999899999999
$ws.Range("H122").Value = 431.72
$ws.Range("I122").Value = 424.5
$ws.Range("J122").Value = 435.11765
$ws.Range("K122").Value = 3820.5
$ws.Range("L122").Value = 3916.05885
$ws.Range("M122").Value = -1370.5
$ws.Range("N122").Value = -8816.058850000001
$ws.Range("H137").Value = 9371.286
$ws.Range("I137").Value = 1049.5
$ws.Range("K137").Value = 3148.5
$ws.Range("M137").Value = 1951.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 68802.42999999999
$ws.Range("J136").Value = 68802.42999999999
$ws.Range("L136").Value = 206407.29
$ws.Range("N136").Value = -211507.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3795.1538
$ws.Range("J7").Value = 3666.6667
$ws.Range("L7").Value = 3666.6667
$ws.Range("N7").Value = -3890.6667
$ws.Range("H46").Value = 3691.4075
$ws.Range("I46").Value = 2350.5
$ws.Range("K46").Value = 2350.5
$ws.Range("M46").Value = -2162.5
$ws.Range("H61").Value = 3828.6191
$ws.Range("I61").Value = 1942.1578
$ws.Range("J61").Value = 21750
$ws.Range("K61").Value = 1942.1578
$ws.Range("L61").Value = 21750
$ws.Range("M61").Value = -1740.1578
$ws.Range("N61").Value = -22154
$ws.Range("H93").Value = 2412.2727
$ws.Range("I93").Value = 2228.5386
$ws.Range("J93").Value = 2677.6667
$ws.Range("K93").Value = 2228.5386
$ws.Range("L93").Value = 2677.6667
$ws.Range("M93").Value = -980.5385999999999
$ws.Range("N93").Value = -5173.6667
$ws.Range("H113").Value = 3828.6191
$ws.Range("I113").Value = 1942.1578
$ws.Range("J113").Value = 21750
$ws.Range("K113").Value = 1942.1578
$ws.Range("L113").Value = 21750
$ws.Range("M113").Value = 227.8422
$ws.Range("N113").Value = -26090
$ws.Range("H126").Value = 3795.1538
$ws.Range("J126").Value = 3666.6667
$ws.Range("L126").Value = 11000.0001
$ws.Range("N126").Value = -15940.0001
$ws.Range("H132").Value = 3100
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("H136").Value = 15177.8
$ws.Range("I136").Value = 13629.667
$ws.Range("J136").Value = 17500
$ws.Range("K136").Value = 40889.001
$ws.Range("L136").Value = 52500
$ws.Range("M136").Value = -38339.001
$ws.Range("N136").Value = -57600
$ws.Range("H139").Value = 67715
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1761.25
$ws.Range("I126").Value = 1739.4546
$ws.Range("K126").Value = 5218.3638
$ws.Range("M126").Value = -2748.3638
$ws.Range("H132").Value = 2925.5334
$ws.Range("I132").Value = 2949
$ws.Range("K132").Value = 8847
$ws.Range("L132").Value = 8847
$ws.Range("M132").Value = -6317
